$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the "Meta description" paragraph that used to sit right
#    under the title heading.
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Meta description*") {
        [void]$p.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------
# 2) Just before the closing "Create a feature image..." paragraph,
#    add a new bold paragraph repeating the page title, and turn the
#    former image-prompt paragraph into the (now italic) meta
#    description text.
#
#    We replace the range spanning the last two paragraphs (the
#    "Lower paying symbols..." bullet and the "Create a feature
#    image..." paragraph) with: the bullet unchanged, the new bold
#    title paragraph, and the italicized meta-description paragraph.
#    Doing the swap across a paragraph boundary (rather than on the
#    very last paragraph alone) avoids leaving a stray empty
#    paragraph behind at the end of the document.
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
$prevPara = $d.Paragraphs($count - 1)
$lastPara = $d.Paragraphs($count)
$r = $d.Range($prevPara.Range.Start, $lastPara.Range.End)

$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>Lower paying symbols have a low payout</w:t></w:r></w:p><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Agent Destiny Free | A Retro Spy-Themed Slot Game</w:t></w:r></w:p><w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Agent Destiny is an exciting spy-themed slot game with a retro comic book style. Play the game for free and enjoy features like Linked Reels and Colossal Symbols.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

[void]$r.InsertXML($xml)
